# Updated cryptos list (price/volume refresh), mirrors the GitHub Actions
# scrape commit. Price cells that look numeric ("522.17", "0.510", ...) are
# forced to text (NumberFormat "@") before assignment so Excel doesn't
# auto-convert them to floating point numbers and lose the exact formatting
# - matching the source data's inline-string representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.515.50'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '3.098.21'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.17'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.097.72'
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("E12").Value = '  +2.22%  '
$ws.Range("D13").Value = '3.632.94'
$ws.Range("E13").Value = '  +1.60%  '
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.66'
$ws.Range("E15").Value = '  -1.98%  '
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '57.614.12'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = '3.104.09'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.09'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '336.86'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.510'
$ws.Range("E24").Value = '  +2.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.46'
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Value = '0.0₃0908'
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("B29").Value = 'USDe'
$ws.Range("C29").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.44'
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.82'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.18'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '155.97'
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.59'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.04'
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0658'
$ws.Range("E40").Value = '  -2.87%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.94'
$ws.Range("E41").Value = '  +0.88%  '
$ws.Range("B42").Value = 'RenzoRestakedETH'
$ws.Range("C42").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D42").Value = '3.139.53'
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("E43").Value = '  +10.92%  '
$ws.Range("E44").Value = '  +4.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.72'
$ws.Range("D47").Value = '2.298.55'
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.977'
$ws.Range("E49").Value = '  +4.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.40'
$ws.Range("E50").Value = '  -0.70%  '
$ws.Range("E51").Value = '  +1.63%  '
